# Add new "Shelled Nut" related error codes to the Error Codes sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18: fill in the previously-empty Meaning cell and update
#     the Class / Function columns for the existing row 16 entry. ---
$ws.Range("C18").Value = "Shelled Nut being loaded but a nut inside it has no type"
$ws.Range("D18").Value = "ShelledNut"
$ws.Range("E18").Value = "Instantiate"

# --- Row 19: brand new error code entry (#17). ---
$ws.Range("B19").Value = 17
$ws.Range("C19").Value = "Shelled Nut being loaded but a nut inside it has no type"
$ws.Range("D19").Value = "Nut"
$ws.Range("E19").Value = "Deserialize"

# Widen column C so the longer "Meaning" text is readable.
$ws.Columns.Item(3).ColumnWidth = 54.5834

# Match the saved selection state left by the author.
$ws.Range("C19").Select()
